# Insert a new row at 246 (this shifts the existing rows 246-356 down to 247-357).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("246:246").Insert()

# The new blank row 246 should start out as a duplicate of the row that is now
# at 247 (the original row 246), then a handful of its cells get overwritten
# with the new values from the commit.
$ws.Range("A247:R247").Copy($ws.Range("A246:R246"))

# Apply the new values for the freshly-inserted row 246.
$ws.Cells.Item(246, 4).Value = 44992           # D246 Fecha
$ws.Cells.Item(246, 10).Value = 50              # J246 Volumen
$ws.Cells.Item(246, 11).Value = 21000           # K246 Precio minimo
$ws.Cells.Item(246, 12).Value = 21000           # L246 Precio maximo
$ws.Cells.Item(246, 13).Value = 21000           # M246 Precio promedio ponderado
$ws.Cells.Item(246, 14).Value = "$/malla 10 kilos"  # N246 Unidad de comercializacion
$ws.Cells.Item(246, 16).Value = 2100            # P246 Precio $/Kg
